$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.212.24"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "'3.153.53"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'591.31"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Value = "'138.46"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'3.152.55"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'5.31"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "'34.30"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "'3.670.71"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "'3.150.28"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'63.205.06"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "'6.68"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").Value = "'478.01"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'14.06"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "'0.704"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'7.74"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "'84.80"
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("D25").Value = "'13.03"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'7.19"
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").Value = "'2.13"
$ws.Range("E30").Value = "  +4.34%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'27.03"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("D34").Value = "'2.55"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.82"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'52.49"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'0.0₃0702"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").Value = "'0.0390"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'418.13"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("D41").Value = "'2.75"
$ws.Range("E41").Value = "  -6.29%  "
$ws.Range("D42").Value = "'8.30"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "'2.931.59"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("E44").Value = "  -6.58%  "
$ws.Range("D45").Value = "'0.264"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.14"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'25.54"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  -6.70%  "
$ws.Range("D51").Value = "'121.67"
$ws.Range("E51").Value = "  +0.16%  "
